# Update cryptos list - values refreshed by GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing it to stay plain text,
# even when the string looks like a decimal number (e.g. "144.33").
# Cells in the "Price" column are plain inline strings with no custom
# number format, so we temporarily mark the cell as Text, assign the
# value, then restore the default "Normal" style (no style index) so
# the cell's XML stays attribute-for-attribute the same as before.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "60.194.27"
$ws.Range("E2").Value = "  +2.21%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.553.66"
$ws.Range("E3").Value = "  +2.62%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  +1.09%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "144.33"
$ws.Range("E6").Value = "  +1.19%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Range("D9") "2.575.79"
$ws.Range("E9").Value = "  +2.24%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.90%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.80%  "

# Row 12 - Toncoin
Set-TextValue $ws.Range("D12") "5.49"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +3.16%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "3.002.73"
$ws.Range("E14").Value = "  +2.36%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "24.03"
$ws.Range("E15").Value = "  +1.19%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "60.177.73"
$ws.Range("E16").Value = "  +2.49%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +3.92%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.560.14"
$ws.Range("E18").Value = "  +1.74%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "11.26"
$ws.Range("E19").Value = "  -1.11%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +1.41%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "327.10"
$ws.Range("E21").Value = "  +1.63%  "

# Row 22 & 23 - Uniswap / Dai swapped positions, with refreshed values
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D22") "1.00"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D23") "5.98"
$ws.Range("E23").Value = "  +4.23%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "63.56"
$ws.Range("E24").Value = "  +4.42%  "

# Row 25 - Polygon
Set-TextValue $ws.Range("D25") "0.433"
$ws.Range("E25").Value = "  -1.06%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +3.99%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  -0.35%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +4.33%  "

# Row 29 - Aptos
$ws.Range("E29").Value = "  +3.26%  "

# Row 30 - PEPE
$ws.Range("E30").Value = "  +4.22%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +2.15%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  -4.05%  "

# Row 33 - Monero
Set-TextValue $ws.Range("D33") "166.27"
$ws.Range("E33").Value = "  +5.82%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  +5.19%  "

# Row 35 - USDe
$ws.Range("E35").Value = "  +0.18%  "

# Row 36 - EthereumClassic
Set-TextValue $ws.Range("D36") "18.74"
$ws.Range("E36").Value = "  +1.25%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  +0.64%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +1.97%  "

# Row 39 - OKB
Set-TextValue $ws.Range("D39") "37.16"
$ws.Range("E39").Value = "  +1.01%  "

# Row 40 - RenderToken
Set-TextValue $ws.Range("D40") "5.61"
$ws.Range("E40").Value = "  -5.58%  "

# Row 41 - Bittensor
Set-TextValue $ws.Range("D41") "301.80"
$ws.Range("E41").Value = "  -2.44%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  +2.14%  "

# Row 43 - SuiNetwork
$ws.Range("E43").Value = "  +5.29%  "

# Row 44 - Mantle
$ws.Range("E44").Value = "  +2.97%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  -0.01%  "

# Row 46 - WhiteBITCoin
Set-TextValue $ws.Range("D46") "10.86"
$ws.Range("E46").Value = "  +0.75%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "127.80"
$ws.Range("E47").Value = "  +2.76%  "

# Row 48 - EnergySwap
Set-TextValue $ws.Range("D48") "19.01"
$ws.Range("E48").Value = "  +2.20%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  +1.68%  "

# Row 50 - Hedera
$ws.Range("E50").Value = "  +0.94%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  +1.21%  "
